$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cycle-time seed value (C2) from 5400 to 2200.
$ws.Range("C2").Value = 2200

# Move the active selection to C2 (matches the saved cursor position in the diff).
$ws.Range("C2").Select()
